# Regenerate the localization-status report for archive:
#  - Files that were "Ready for handoff" have moved on to "In Translation".
#  - The status columns are narrower now that the new status text is shorter,
#    so re-fit/resize them to match.

$wb = $excel.ActiveWorkbook

# --- 1. Update status text wherever it currently reads "Ready for handoff" ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if ($v -is [string] -and $v -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Resize the status columns to fit the new, shorter text ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (Status)
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
